# Updates cryptos list cell values to match the latest scrape.
# Uses a leading apostrophe to force text interpretation for numeric-looking
# values (keeping cell type = Text, matching the original inlineStr cells),
# then resets the cell style to "Normal" so no stray number-format style is
# left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'68.854.22"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "'  -1.46%  "
$ws.Cells.Item(2, 5).Style = "Normal"

$ws.Cells.Item(3, 4).Value = "'3.425.27"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "'  -2.09%  "
$ws.Cells.Item(3, 5).Style = "Normal"

$ws.Cells.Item(4, 4).Value = "'0.998"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "'  -0.28%  "
$ws.Cells.Item(4, 5).Style = "Normal"

$ws.Cells.Item(5, 4).Value = "'609.81"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "'  +1.58%  "
$ws.Cells.Item(5, 5).Style = "Normal"

$ws.Cells.Item(6, 4).Value = "'166.03"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "'  -3.95%  "
$ws.Cells.Item(6, 5).Style = "Normal"

$ws.Cells.Item(7, 4).Value = "'3.423.43"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "'  -1.87%  "
$ws.Cells.Item(7, 5).Style = "Normal"

$ws.Cells.Item(8, 5).Value = "'  -2.41%  "
$ws.Cells.Item(8, 5).Style = "Normal"

$ws.Cells.Item(9, 5).Value = "'  +0.14%  "
$ws.Cells.Item(9, 5).Style = "Normal"

$ws.Cells.Item(10, 4).Value = "'0.192"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "'  +0.04%  "
$ws.Cells.Item(10, 5).Style = "Normal"

$ws.Cells.Item(11, 5).Value = "'  -4.11%  "
$ws.Cells.Item(11, 5).Style = "Normal"

$ws.Cells.Item(12, 4).Value = "'0.559"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "'  -3.39%  "
$ws.Cells.Item(12, 5).Style = "Normal"

$ws.Cells.Item(13, 4).Value = "'43.72"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "'  -4.48%  "
$ws.Cells.Item(13, 5).Style = "Normal"

$ws.Cells.Item(14, 4).Value = "'0.0000267"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "'  -2.35%  "
$ws.Cells.Item(14, 5).Style = "Normal"

$ws.Cells.Item(15, 4).Value = "'3.975.22"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "'  -2.19%  "
$ws.Cells.Item(15, 5).Style = "Normal"

$ws.Cells.Item(16, 4).Value = "'8.07"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "'  -1.93%  "
$ws.Cells.Item(16, 5).Style = "Normal"

$ws.Cells.Item(17, 2).Value = "WrappedBTC"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(17, 4).Value = "'68.924.42"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "'  -1.48%  "
$ws.Cells.Item(17, 5).Style = "Normal"

$ws.Cells.Item(18, 2).Value = "WrappedEther"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(18, 4).Value = "'3.420.93"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "'  -2.46%  "
$ws.Cells.Item(18, 5).Style = "Normal"

$ws.Cells.Item(19, 4).Value = "'576.32"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "'  -5.42%  "
$ws.Cells.Item(19, 5).Style = "Normal"

$ws.Cells.Item(20, 5).Value = "'  +0.56%  "
$ws.Cells.Item(20, 5).Style = "Normal"

$ws.Cells.Item(21, 4).Value = "'17.01"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "'  -0.70%  "
$ws.Cells.Item(21, 5).Style = "Normal"

$ws.Cells.Item(22, 4).Value = "'0.839"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "'  -3.19%  "
$ws.Cells.Item(22, 5).Style = "Normal"

$ws.Cells.Item(23, 4).Value = "'8.83"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "'  -0.48%  "
$ws.Cells.Item(23, 5).Style = "Normal"

$ws.Cells.Item(24, 4).Value = "'95.14"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "'  -0.54%  "
$ws.Cells.Item(24, 5).Style = "Normal"

$ws.Cells.Item(25, 4).Value = "'14.99"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "'  -3.08%  "
$ws.Cells.Item(25, 5).Style = "Normal"

$ws.Cells.Item(26, 4).Value = "'3.60"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "'  -2.70%  "
$ws.Cells.Item(26, 5).Style = "Normal"

$ws.Cells.Item(27, 5).Value = "'  +0.12%  "
$ws.Cells.Item(27, 5).Style = "Normal"

$ws.Cells.Item(28, 4).Value = "'2.39"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "'  -5.90%  "
$ws.Cells.Item(28, 5).Style = "Normal"

$ws.Cells.Item(29, 4).Value = "'32.49"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "'  -4.47%  "
$ws.Cells.Item(29, 5).Style = "Normal"

$ws.Cells.Item(30, 4).Value = "'8.49"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "'  -4.05%  "
$ws.Cells.Item(30, 5).Style = "Normal"

$ws.Cells.Item(31, 4).Value = "'7.74"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "'  -3.70%  "
$ws.Cells.Item(31, 5).Style = "Normal"

$ws.Cells.Item(32, 5).Value = "'  -3.79%  "
$ws.Cells.Item(32, 5).Style = "Normal"

$ws.Cells.Item(33, 5).Value = "'  -7.04%  "
$ws.Cells.Item(33, 5).Style = "Normal"

$ws.Cells.Item(34, 4).Value = "'6.50"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "'  -4.21%  "
$ws.Cells.Item(34, 5).Style = "Normal"

$ws.Cells.Item(35, 4).Value = "'581.23"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "'  -6.84%  "
$ws.Cells.Item(35, 5).Style = "Normal"

$ws.Cells.Item(36, 4).Value = "'10.39"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "'  -2.38%  "
$ws.Cells.Item(36, 5).Style = "Normal"

$ws.Cells.Item(37, 5).Value = "'  +0.23%  "
$ws.Cells.Item(37, 5).Style = "Normal"

$ws.Cells.Item(38, 2).Value = "Hedera"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(38, 4).Value = "'0.0947"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "'  -4.05%  "
$ws.Cells.Item(38, 5).Style = "Normal"

$ws.Cells.Item(39, 2).Value = "VeChain"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(39, 4).Value = "'0.0464"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "'  -2.80%  "
$ws.Cells.Item(39, 5).Style = "Normal"

$ws.Cells.Item(40, 4).Value = "'55.77"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "'  -0.59%  "
$ws.Cells.Item(40, 5).Style = "Normal"

$ws.Cells.Item(41, 5).Value = "'  -0.63%  "
$ws.Cells.Item(41, 5).Style = "Normal"

$ws.Cells.Item(42, 4).Value = "'3.07"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "'  -13.48%  "
$ws.Cells.Item(42, 5).Style = "Normal"

$ws.Cells.Item(43, 4).Value = "'3.220.16"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "'  -3.11%  "
$ws.Cells.Item(43, 5).Style = "Normal"

$ws.Cells.Item(44, 4).Value = "'0.0₃0671"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "'  -6.44%  "
$ws.Cells.Item(44, 5).Style = "Normal"

$ws.Cells.Item(45, 4).Value = "'30.80"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "'  -3.56%  "
$ws.Cells.Item(45, 5).Style = "Normal"

$ws.Cells.Item(46, 5).Value = "'  -5.55%  "
$ws.Cells.Item(46, 5).Style = "Normal"

$ws.Cells.Item(47, 5).Value = "'  -4.79%  "
$ws.Cells.Item(47, 5).Style = "Normal"

$ws.Cells.Item(48, 5).Value = "'  -6.44%  "
$ws.Cells.Item(48, 5).Style = "Normal"

$ws.Cells.Item(49, 5).Value = "'  -3.00%  "
$ws.Cells.Item(49, 5).Style = "Normal"

$ws.Cells.Item(50, 4).Value = "'132.01"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "'  -1.40%  "
$ws.Cells.Item(50, 5).Style = "Normal"

$ws.Cells.Item(51, 5).Value = "'  +0.01%  "
$ws.Cells.Item(51, 5).Style = "Normal"
